$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.140973091125488
$ws.Range("B1").Value = 2.328085660934448
$ws.Range("C1").Value = 4.846299648284912
$ws.Range("D1").Value = 2.182245492935181
$ws.Range("E1").Value = 1.075821757316589
